# "First cut of better Estimands"
# Reorders the treatmentXref / endpointXref / intercurrentEventStrategy
# columns on the studyDesignEstimands sheet and appends new estimand rows,
# then leaves that sheet active/selected with a resized workbook window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studyDesignEstimands")

# --- Reorder columns F/G/H -------------------------------------------------
# Before: F=intercurrentEventStrategy, G=treatmentXref,           H=endpointXref
# After:  F=treatmentXref,             G=endpointXref,            H=intercurrentEventStrategy
$strategyHeader = $ws.Cells.Item(1, 6).Value2
$treatmentXrefHeader = $ws.Cells.Item(1, 7).Value2
$endpointXrefHeader = $ws.Cells.Item(1, 8).Value2

$strategyValue = $ws.Cells.Item(2, 6).Value2
$treatmentXrefValue = $ws.Cells.Item(2, 7).Value2
$endpointXrefValue = $ws.Cells.Item(2, 8).Value2

$ws.Cells.Item(1, 6).Value = $treatmentXrefHeader
$ws.Cells.Item(1, 7).Value = $endpointXrefHeader
$ws.Cells.Item(1, 8).Value = $strategyHeader

$ws.Cells.Item(2, 6).Value = $treatmentXrefValue
$ws.Cells.Item(2, 7).Value = $endpointXrefValue
$ws.Cells.Item(2, 8).Value = $strategyValue

# Column widths follow the moved content: F/G become the narrow xref
# columns (17 chars) and H becomes the wide free-text column (61.83 chars).
# (ColumnWidth is specified in characters; COM rounds to whole display
# pixels like real Excel does, so feed it the inverse of that rounding.)
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 16.166666666666668
$ws.Columns.Item(8).ColumnWidth = 60.998697916666664

# --- New intercurrence-event free text rows (H3:H4) -------------------------
$ws.Cells.Item(3, 8).Value = "A second bad event"
$ws.Cells.Item(4, 8).Value = "A third bad thing"

# --- New estimand row (row 5) ----------------------------------------------
$ws.Cells.Item(5, 1).Value = "EST2"
$ws.Cells.Item(5, 2).Value = "Something else"
$ws.Cells.Item(5, 3).Value = "ITT"
$ws.Cells.Item(5, 4).Value = "Bad stuff"
$ws.Cells.Item(5, 5).Value = "IC Event Description Number 2"
$ws.Cells.Item(5, 6).Value = "INT2"
$ws.Cells.Item(5, 7).Value = "END2"
$ws.Cells.Item(5, 8).Value = "Really really bad shit"

# --- Selection / active sheet -----------------------------------------------
$ws.Activate()
$ws.Range("D29").Select()
